$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row in the sheet
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 4)  # Column D = Speaker
    if ($cell.Value2 -eq "HILLARY LEWIS-WOLFSEN") {
        $cell.Value = "T"
    }
}

$wb.Save()
